$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shifted rows 1109-1201: columns D, I, J, K, L, M, O, P
$rowsData = @(
    @{Row=1109; D=45223; I="Primera"; J=2200; K=700; L=750; M=723; O="Provincia de Quillota"; P=723}
    @{Row=1110; D=45223; I="Segunda"; J=850; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1111; D=44468; I="Primera"; J=2200; K=650; L=700; M=677; O="Provincia de Quillota"; P=677}
    @{Row=1112; D=44468; I="Segunda"; J=1100; K=550; L=550; M=550; O="Provincia de Quillota"; P=550}
    @{Row=1113; D=44487; I="Primera"; J=2650; K=600; L=700; M=653; O="Provincia de Quillota"; P=653}
    @{Row=1114; D=44487; I="Segunda"; J=1600; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1115; D=44714; I="Primera"; J=2400; K=1000; L=1100; M=1050; O="Provincia de Quillota"; P=1050}
    @{Row=1116; D=44714; I="Segunda"; J=1100; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1117; D=44778; I="Primera"; J=2400; K=1000; L=1100; M=1050; O="Provincia de Quillota"; P=1050}
    @{Row=1118; D=44778; I="Segunda"; J=1150; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1119; D=44754; I="Primera"; J=1500; K=1100; L=1100; M=1100; O="Provincia de Quillota"; P=1100}
    @{Row=1120; D=44754; I="Segunda"; J=1200; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1121; D=44462; I="Primera"; J=3300; K=600; L=650; M=624; O="Provincia de Quillota"; P=624}
    @{Row=1122; D=44462; I="Segunda"; J=1500; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1123; D=44490; I="Primera"; J=2700; K=600; L=650; M=624; O="Provincia de Quillota"; P=624}
    @{Row=1124; D=44490; I="Segunda"; J=1500; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1125; D=45215; I="Primera"; J=2800; K=700; L=800; M=743; O="Provincia de Quillota"; P=743}
    @{Row=1126; D=44238; I="Primera"; J=1400; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000}
    @{Row=1127; D=44238; I="Segunda"; J=1200; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1128; D=44924; I="Primera"; J=1800; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1129; D=44924; I="Segunda"; J=1100; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1130; D=44973; I="Primera"; J=1500; K=1000; L=1100; M=1050; O="Provincia de Quillota"; P=1050}
    @{Row=1131; D=44973; I="Segunda"; J=780; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1132; D=44874; I="Primera"; J=2400; K=750; L=800; M=775; O="Provincia de Quillota"; P=775}
    @{Row=1133; D=44874; I="Segunda"; J=1300; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1134; D=44245; I="Primera"; J=1100; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1135; D=44245; I="Segunda"; J=950; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1136; D=44481; I="Primera"; J=850; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1137; D=44481; I="Segunda"; J=970; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1138; D=44999; I="Primera"; J=2400; K=1000; L=1100; M=1050; O="Provincia de Quillota"; P=1050}
    @{Row=1139; D=44999; I="Segunda"; J=1388; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1140; D=44294; I="Primera"; J=2500; K=850; L=900; M=874; O="Provincia de Quillota"; P=874}
    @{Row=1141; D=44294; I="Segunda"; J=1300; K=750; L=750; M=750; O="Provincia de Quillota"; P=750}
    @{Row=1142; D=44558; I="Primera"; J=2500; K=600; L=650; M=626; O="Provincia de Quillota"; P=626}
    @{Row=1143; D=44558; I="Segunda"; J=1200; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1144; D=44634; I="Primera"; J=930; K=1100; L=1200; M=1152; O="Provincia de Quillota"; P=1152}
    @{Row=1145; D=44634; I="Segunda"; J=460; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1146; D=44174; I="Primera"; J=900; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1147; D=44174; I="Segunda"; J=958; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1148; D=44608; I="Primera"; J=880; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000}
    @{Row=1149; D=44608; I="Segunda"; J=850; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1150; D=44550; I="Primera"; J=1800; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1151; D=44550; I="Segunda"; J=1900; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1152; D=44775; I="Primera"; J=2550; K=1100; L=1200; M=1157; O="Provincia de Quillota"; P=1157}
    @{Row=1153; D=44775; I="Segunda"; J=1200; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1154; D=44859; I="Primera"; J=4900; K=850; L=900; M=873; O="Provincia de Quillota"; P=873}
    @{Row=1155; D=44859; I="Segunda"; J=1900; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1156; D=44910; I="Primera"; J=1800; K=850; L=900; M=876; O="Provincia de Quillota"; P=876}
    @{Row=1157; D=44910; I="Segunda"; J=900; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1158; D=45063; I="Primera"; J=3100; K=1100; L=1200; M=1152; O="Provincia de Quillota"; P=1152}
    @{Row=1159; D=45063; I="Segunda"; J=1400; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1160; D=44365; I="Primera"; J=2700; K=700; L=750; M=717; O="Provincia de Quillota"; P=717}
    @{Row=1161; D=44365; I="Segunda"; J=1600; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1162; D=45069; I="Primera"; J=3400; K=1100; L=1200; M=1153; O="Provincia de Quillota"; P=1153}
    @{Row=1163; D=45069; I="Segunda"; J=1700; K=1000; L=1000; M=1000; O="Provincia de Quillota"; P=1000}
    @{Row=1164; D=44711; I="Primera"; J=3100; K=850; L=900; M=874; O="Provincia de Quillota"; P=874}
    @{Row=1165; D=44711; I="Segunda"; J=1200; K=650; L=650; M=650; O="Provincia de Quillota"; P=650}
    @{Row=1166; D=44382; I="Primera"; J=1900; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1167; D=44382; I="Segunda"; J=1850; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1168; D=44795; I="Primera"; J=2500; K=1100; L=1200; M=1148; O="Provincia de Quillota"; P=1148}
    @{Row=1169; D=44795; I="Segunda"; J=1200; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1170; D=45173; I="Primera"; J=3200; K=700; L=750; M=725; O="Provincia de Quillota"; P=725}
    @{Row=1171; D=44883; I="Primera"; J=1600; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1172; D=44883; I="Segunda"; J=1200; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1173; D=44673; I="Primera"; J=1200; K=1100; L=1100; M=1100; O="Provincia de Quillota"; P=1100}
    @{Row=1174; D=44673; I="Segunda"; J=1100; K=850; L=850; M=850; O="Provincia de Quillota"; P=850}
    @{Row=1175; D=44818; I="Primera"; J=3700; K=800; L=1100; M=962; O="Provincia de Quillota"; P=962}
    @{Row=1176; D=44649; I="Primera"; J=1300; K=1200; L=1300; M=1250; O="Provincia de Quillota"; P=1250}
    @{Row=1177; D=44649; I="Segunda"; J=550; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1178; D=44341; I="Primera"; J=1800; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1179; D=44341; I="Segunda"; J=1500; K=700; L=700; M=700; O="Provincia de Quillota"; P=700}
    @{Row=1180; D=44868; I="Primera"; J=3200; K=900; L=950; M=925; O="Provincia de Quillota"; P=925}
    @{Row=1181; D=44868; I="Segunda"; J=2400; K=700; L=750; M=719; O="Provincia de Quillota"; P=719}
    @{Row=1182; D=45216; I="Primera"; J=3200; K=700; L=800; M=750; O="Provincia de Quillota"; P=750}
    @{Row=1183; D=45216; I="Segunda"; J=1000; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1184; D=44980; I="Primera"; J=3100; K=1000; L=1100; M=1052; O="Provincia de Quillota"; P=1052}
    @{Row=1185; D=44980; I="Segunda"; J=1200; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1186; D=44460; I="Primera"; J=2300; K=600; L=650; M=626; O="Provincia de Quillota"; P=626}
    @{Row=1187; D=45114; I="Primera"; J=3400; K=700; L=800; M=747; O="Provincia de Quillota"; P=747}
    @{Row=1188; D=45114; I="Segunda"; J=1200; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1189; D=45012; I="Primera"; J=2200; K=1100; L=1200; M=1150; O="Provincia de Quillota"; P=1150}
    @{Row=1190; D=45012; I="Segunda"; J=1200; K=900; L=900; M=900; O="Provincia de Quillota"; P=900}
    @{Row=1191; D=44232; I="Primera"; J=1800; K=800; L=850; M=825; O="Provincia de Copiapó"; P=825}
    @{Row=1192; D=44232; I="Segunda"; J=880; K=650; L=650; M=650; O="Provincia de Copiapó"; P=650}
    @{Row=1193; D=44613; I="Primera"; J=550; K=1100; L=1100; M=1100; O="Provincia de Melipilla"; P=1100}
    @{Row=1194; D=44613; I="Segunda"; J=450; K=800; L=800; M=800; O="Provincia de Melipilla"; P=800}
    @{Row=1195; D=44725; I="Primera"; J=2250; K=1000; L=1100; M=1049; O="Provincia de Quillota"; P=1049}
    @{Row=1196; D=44725; I="Segunda"; J=1250; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1197; D=44893; I="Primera"; J=920; K=800; L=800; M=800; O="Provincia de Quillota"; P=800}
    @{Row=1198; D=44893; I="Segunda"; J=950; K=650; L=650; M=650; O="Provincia de Quillota"; P=650}
    @{Row=1199; D=44392; I="Primera"; J=2600; K=650; L=700; M=679; O="Provincia de Quillota"; P=679}
    @{Row=1200; D=44565; I="Primera"; J=3300; K=600; L=650; M=627; O="Provincia de Quillota"; P=627}
    @{Row=1201; D=44565; I="Segunda"; J=1200; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
    @{Row=1202; D=44544; I="Primera"; J=1600; K=600; L=600; M=600; O="Provincia de Quillota"; P=600}
    @{Row=1203; D=44544; I="Segunda"; J=1500; K=500; L=500; M=500; O="Provincia de Quillota"; P=500}
)

foreach ($item in $rowsData) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
}

# New rows 1202 and 1203: set all columns A-R
$ws.Cells.Item(1202, 1).Value = 3
$ws.Cells.Item(1202, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1202, 3).Value = "Coquimbo"
$ws.Cells.Item(1202, 5).Value = 5
$ws.Cells.Item(1202, 6).Value = 100112008
$ws.Cells.Item(1202, 7).Value = "Coliflor"
$ws.Cells.Item(1202, 8).Value = "Sin especificar"
$ws.Cells.Item(1202, 14).Value = "`$/unidad"
$ws.Cells.Item(1202, 17).Value = 1
$ws.Cells.Item(1202, 18).Value = "Hortaliza"
$ws.Cells.Item(1203, 1).Value = 3
$ws.Cells.Item(1203, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1203, 3).Value = "Coquimbo"
$ws.Cells.Item(1203, 5).Value = 5
$ws.Cells.Item(1203, 6).Value = 100112008
$ws.Cells.Item(1203, 7).Value = "Coliflor"
$ws.Cells.Item(1203, 8).Value = "Sin especificar"
$ws.Cells.Item(1203, 14).Value = "`$/unidad"
$ws.Cells.Item(1203, 17).Value = 1
$ws.Cells.Item(1203, 18).Value = "Hortaliza"
